# Update cryptos list (Price / Volume(1h) columns) as produced by the
# scheduled GitHub Actions refresh on Sat Apr  6 17:38:54 UTC 2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new Price, new Volume(1h)). $null means "leave unchanged".
$updates = @{
    2  = @("68.115.58", "  +0.61%  ")
    3  = @("3.338.62",  "  +0.39%  ")
    4  = @("0.999",     "  +0.07%  ")
    5  = @("584.52",    "  +0.42%  ")
    6  = @("177.05",    $null)
    7  = @("0.999",     "  -0.09%  ")
    8  = @($null,       "  +1.38%  ")
    10 = @($null,       "  +1.45%  ")
    11 = @($null,       "  +5.57%  ")
    12 = @($null,       "  +2.12%  ")
    13 = @("694.78",    "  +3.69%  ")
    14 = @("3.882.27",  "  +0.58%  ")
    15 = @($null,       "  +0.75%  ")
    16 = @("68.151.87", "  +0.46%  ")
    17 = @($null,       "  +1.32%  ")
    18 = @("3.335.42",  "  +0.35%  ")
    19 = @("17.49",     "  +0.39%  ")
    20 = @("11.17",     "  +2.58%  ")
    21 = @($null,       "  +0.86%  ")
    22 = @("5.43",      "  +0.73%  ")
    23 = @("16.97",     "  -0.04%  ")
    24 = @("100.62",    "  +3.03%  ")
    25 = @($null,       "  +1.85%  ")
    26 = @($null,       "  +0.96%  ")
    27 = @("9.50",      "  +2.10%  ")
    28 = @("33.14",     "  -0.85%  ")
    29 = @($null,       "  +1.72%  ")
    30 = @("6.97",      "  -4.87%  ")
    31 = @("567.93",    "  -3.06%  ")
    32 = @($null,       "  +1.13%  ")
    33 = @($null,       "  +1.66%  ")
    34 = @("3.730.94",  "  +0.36%  ")
    35 = @($null,       "  +1.10%  ")
    36 = @("0.998",     "  -0.15%  ")
    37 = @($null,       "  +2.82%  ")
    38 = @($null,       "  +3.60%  ")
    39 = @("35.14",     "  +7.27%  ")
    40 = @("3.17",      "  +2.46%  ")
    41 = @($null,       "  +0.06%  ")
    42 = @($null,       "  +1.94%  ")
    43 = @($null,       "  +0.96%  ")
    44 = @($null,       "  +0.49%  ")
    45 = @($null,       "  +1.43%  ")
    46 = @("2.65",      "  +2.59%  ")
    47 = @($null,       "  +1.23%  ")
    48 = @($null,       "  -0.05%  ")
    49 = @($null,       "  -0.27%  ")
    50 = @("130.72",    "  +3.07%  ")
    51 = @($null,       "  +0.64%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $priceVal = $vals[0]
    $volVal = $vals[1]

    # Some "prices" (e.g. "0.999", "17.49") parse as plain numbers, so a
    # naive assignment would silently convert the cell to a numeric type.
    # The source sheet always stores these as text, so force text storage
    # via a temporary "@" number format, then restore the default ("Normal")
    # style afterwards so no stray style index is left on the cell.
    if ($null -ne $priceVal) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Formula = $priceVal
        $cell.Style = "Normal"
    }
    if ($null -ne $volVal) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.NumberFormat = "@"
        $cell.Formula = $volVal
        $cell.Style = "Normal"
    }
}
